$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting used by the other header cells (e.g. H1) onto the new headers
$ws.Cells.Item(1, 8).Copy()
$ws.Cells.Item(1, 9).PasteSpecial(-4122)
$ws.Cells.Item(1, 10).PasteSpecial(-4122)

# New header cells for columns I and J
$ws.Cells.Item(1, 9).Value = "I0"
$ws.Cells.Item(1, 10).Value = "IF"

# Data for columns I (I0) and J (IF), rows 2-28
$data = @(
    @{Row=2; I=1; J=3},
    @{Row=3; I=1; J=5},
    @{Row=4; I=1; J=5},
    @{Row=5; I=1; J=5},
    @{Row=6; I=1; J=4},
    @{Row=7; I=1; J=6},
    @{Row=8; I=1; J=6},
    @{Row=9; I=1; J=5},
    @{Row=10; I=1; J=6},
    @{Row=11; I=1; J=6},
    @{Row=12; I=1; J=7},
    @{Row=13; I=1; J=5},
    @{Row=14; I=1; J=4},
    @{Row=15; I=1; J=4},
    @{Row=16; I=1; J=3},
    @{Row=17; I=1; J=6},
    @{Row=18; I=1; J=6},
    @{Row=19; I=1; J=7},
    @{Row=20; I=1; J=7},
    @{Row=21; I=1; J=6},
    @{Row=22; I=1; J=5},
    @{Row=23; I=1; J=6},
    @{Row=24; I=1; J=5},
    @{Row=25; I=1; J=6},
    @{Row=26; I=9; J=9},
    @{Row=27; I=6; J=7},
    @{Row=28; I=1; J=1}
)

foreach ($entry in $data) {
    $ws.Cells.Item($entry.Row, 9).Value = $entry.I
    $ws.Cells.Item($entry.Row, 10).Value = $entry.J
}
